# Insn-ASCII conversion - add remaining instructions to the main loop table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows -------------------------------------------------
# Row 9: "UND" was a placeholder/typo - it's actually "UNDO"
$ws.Range("B9").Value = "UNDO"

# Row 10: "RDO" -> "REDO", with corresponding hex/decimal/binary updates
$ws.Range("B10").Value = "REDO"
$ws.Range("C10").Value = "524544"
$ws.Range("D10").Value = 5391684
$ws.Range("E10").Value = "010100100100010101000100"

# --- Fill down the SQRT/ROUND/square/remainder helper formulas ---------
# Row 12 (blank spacer row) now also carries the helper formulas
$ws.Range("F12").Formula = "=SQRT(D12)"
$ws.Range("G12").Formula = "=ROUND(F12, 0)"
$ws.Range("H12").Formula = "=G12*G12"
$ws.Range("I12").Formula = "=D12-H12"
$ws.Range("I12").NumberFormat = "@"

# Row 13 (CLC)
$ws.Range("F13").Formula = "=SQRT(D13)"
$ws.Range("G13").Formula = "=ROUND(F13, 0)"
$ws.Range("H13").Formula = "=G13*G13"
$ws.Range("I13").Formula = "=D13-H13"
$ws.Range("I13").NumberFormat = "@"

# Row 14 (RPT)
$ws.Range("F14").Formula = "=SQRT(D14)"
$ws.Range("G14").Formula = "=ROUND(F14, 0)"
$ws.Range("H14").Formula = "=G14*G14"
$ws.Range("I14").Formula = "=D14-H14"
$ws.Range("I14").NumberFormat = "@"

# --- New instructions: PNUP / PNDN --------------------------------------
# Row 15: PNUP
$ws.Range("B15").Value = "PNUP"
$ws.Range("C15").Value = "50 4e 55"
$ws.Range("D15").Value = 5262933
$ws.Range("E15").Value = "01010000 01001110 01010101"
$ws.Range("F15").Formula = "=SQRT(D15)"
$ws.Range("G15").Formula = "=ROUND(F15, 0)"
$ws.Range("H15").Formula = "=G15*G15"
$ws.Range("I15").Formula = "=D15-H15"
$ws.Range("I15").NumberFormat = "@"

# Row 16: PNDN
$ws.Range("B16").Value = "PNDN"
$ws.Range("C16").Value = "50 4e 44"
$ws.Range("D16").Value = 5262916
$ws.Range("E16").Value = "01010000 01001110 01000100"
$ws.Range("F16").Formula = "=SQRT(D16)"
$ws.Range("G16").Formula = "=ROUND(F16, 0)"
$ws.Range("H16").Formula = "=G16*G16"
$ws.Range("I16").Formula = "=D16-H16"
$ws.Range("I16").NumberFormat = "@"

# Move active selection to B11, matching where editing left off
[void]$ws.Range("B11").Select()
